$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4" (same column
#    layout/headers/styles) and placing it right after "2021-Q4" (i.e. right
#    before "总计").
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newWs = $wb.Worksheets.Item("2021-Q4 (2)")
$newWs.Name = "2022-Q1"

# The duplicated sheet has 3 data rows (rows 2-4); we need 4 data rows
# (rows 2-5), so extend the formatting of the A column down to row 5 by
# copying the format from A4 (which already carries the correct style).
$newWs.Range("A4").Copy()
$newWs.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the fund-holding data for 2022-Q1.
function Set-TextValue($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2
$newWs.Range("A2").Value = 0
Set-TextValue $newWs "B2" "012744"
Set-TextValue $newWs "C2" "光大保德信品质生活混合型证券投资基金A"
Set-TextValue $newWs "D2" "6.91"
Set-TextValue $newWs "E2" "84.96"
Set-TextValue $newWs "F2" "6.74"
Set-TextValue $newWs "G2" "0.4657"
$newWs.Range("H2").Value = 4

# Row 3
$newWs.Range("A3").Value = 1
Set-TextValue $newWs "B3" "007592"
Set-TextValue $newWs "C3" "华夏价值精选混合"
Set-TextValue $newWs "D3" "2.55"
Set-TextValue $newWs "E3" "94.58"
Set-TextValue $newWs "F3" "7.59"
Set-TextValue $newWs "G3" "0.1935"
$newWs.Range("H3").Value = 2

# Row 4
$newWs.Range("A4").Value = 2
Set-TextValue $newWs "B4" "004099"
Set-TextValue $newWs "C4" "前海开源沪港深景气行业精选灵活配置混合"
Set-TextValue $newWs "D4" "0.41"
Set-TextValue $newWs "E4" "93.07"
Set-TextValue $newWs "F4" "7.67"
Set-TextValue $newWs "G4" "0.0314"
$newWs.Range("H4").Value = 10

# Row 5 (new row)
$newWs.Range("A5").Value = 3
Set-TextValue $newWs "B5" "012758"
Set-TextValue $newWs "C5" "光大保德信品质生活混合型证券投资基金C"
Set-TextValue $newWs "D5" "0.31"
Set-TextValue $newWs "E5" "84.96"
Set-TextValue $newWs "F5" "6.74"
Set-TextValue $newWs "G5" "0.0209"
$newWs.Range("H5").Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" (Total) summary sheet: insert a new row for 2022-Q1
#    at the top of the data, push the existing rows down, and renumber the
#    index column (A).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("A2:D2").Style = "Normal"

# Give the new A2 cell the same formatting as the rest of the index column.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.71

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
